$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "col A ColumnWidth:" $ws.Columns.Item(1).ColumnWidth
Write-Host "col A Width:" $ws.Columns.Item(1).Width
Write-Host "col F ColumnWidth:" $ws.Columns.Item(6).ColumnWidth
Write-Host "col F Width:" $ws.Columns.Item(6).Width
Write-Host "Range A1 ColumnWidth:" $ws.Range("A1").ColumnWidth
Write-Host "Range F1 ColumnWidth:" $ws.Range("F1").ColumnWidth
Write-Host "standard width:" $ws.StandardWidth
